$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the selected/active cell shown in the saved view
$ws.Range("H24").Select()

# Fill in row 18 with the final measurement data (100% row of the table)
$ws.Range("B18").Value = 0.0000452561314841659
$ws.Range("C18").Value = 0.000038013758665329
$ws.Range("D18").Value = 0.989694550037384
$ws.Range("E18").Value = 0.989787264585495
$ws.Range("F18").Value = 45.5007385253906
$ws.Range("G18").Value = 45.5530184326171
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = 0.0000301593978929304
$ws.Range("J18").Value = 0.0000385957620019326
$ws.Range("K18").Value = 0.989657027721405
$ws.Range("L18").Value = 0.988322944641113
$ws.Range("M18").Value = 46.9347955703735
$ws.Range("N18").Value = 45.7758985595703
